{"js": "// Update the date heading and the 25 division problems in the practice table.\n// The document has one leading paragraph with the date, followed by a single\n// table whose rows 0, 4, 8, 12, 16 (0-based) hold the 5 visible problems per\n// row (the rows in between are intentionally blank \"answer\" rows).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// 1) Update the date paragraph (first paragraph in the document).\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.load(\"text\");\nawait context.sync();\nif (dateParagraph.text.trim() === \"2026-02-16 Monday\") {\n  const dateRange = dateParagraph.getRange();\n  dateRange.insertText(\"2026-02-17 Tuesday\", \"Replace\");\n}\n\n// 2) Update the division problems inside the table, cell by cell, in document order.\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Map of rowIndex -> [oldText -> newText] pairs, in left-to-right cell order.\nconst rowEdits = {\n  0: [\"51\u00f74=\", \"42\u00f76=\", \"20\u00f73=\", \"83\u00f72=\", \"15\u00f73=\", \"98\u00f75=\", \"10\u00f78=\", \"48\u00f75=\", \"88\u00f78=\", \"61\u00f77=\"],\n  4: [\"34\u00f79=\", \"71\u00f79=\", \"14\u00f78=\", \"66\u00f76=\", \"48\u00f74=\", \"91\u00f72=\", \"40\u00f73=\", \"20\u00f77=\", \"31\u00f78=\", \"59\u00f73=\"],\n  8: [\"98\u00f73=\", \"95\u00f74=\", \"30\u00f74=\", \"73\u00f76=\", \"31\u00f76=\", \"36\u00f77=\", \"97\u00f79=\", \"56\u00f78=\", \"13\u00f73=\", \"28\u00f77=\"],\n  12: [\"82\u00f74=\", \"20\u00f77=\", \"24\u00f74=\", \"84\u00f76=\", \"64\u00f72=\", \"48\u00f74=\", \"33\u00f78=\", \"53\u00f76=\", \"88\u00f75=\", \"27\u00f74=\"],\n  16: [\"72\u00f74=\", \"17\u00f77=\", \"15\u00f73=\", \"76\u00f79=\", \"67\u00f76=\", \"29\u00f73=\", \"37\u00f77=\", \"91\u00f75=\", \"60\u00f72=\", \"59\u00f75=\"],\n};\n\nfor (const rowIndexStr of Object.keys(rowEdits)) {\n  const rowIndex = Number(rowIndexStr);\n  const pairs = rowEdits[rowIndex];\n  const row = rows.items[rowIndex];\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n\n  for (let c = 0; c < cells.items.length; c++) {\n    const oldText = pairs[c * 2];\n    const newText = pairs[c * 2 + 1];\n    const cellParagraphs = cells.items[c].body.paragraphs;\n    cellParagraphs.load(\"items/text\");\n    await context.sync();\n    const cellParagraph = cellParagraphs.items[0];\n    if (cellParagraph.text.trim() === oldText) {\n      // Replace the text on the paragraph's own range so the existing run\n      // (font/size) and paragraph (alignment) formatting is preserved,\n      // instead of rewriting the whole cell body (which would drop it).\n      cellParagraph.getRange().insertText(newText, \"Replace\");\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the 25 division problems in the practice table.\n# The document has one leading paragraph with the date, followed by a single\n# table whose rows 1, 5, 9, 13, 17 (1-based) hold the 5 visible problems per\n# row (the rows in between are intentionally blank \"answer\" rows).\n\n$d = $word.ActiveDocument\n\n# 1) Update the date paragraph (first paragraph in the document).\n$dateParagraph = $d.Paragraphs.Item(1)\n$dateText = $dateParagraph.Range.Text.TrimEnd([char]13, [char]7)\nif ($dateText -eq \"2026-02-16 Monday\") {\n    $dateParagraph.Range.Text = \"2026-02-17 Tuesday\"\n}\n\n# 2) Update the division problems inside the table, cell by cell, in document order.\n$table = $d.Tables.Item(1)\n\n# Map of rowIndex (1-based) -> [oldText, newText] pairs, in left-to-right cell order.\n$rowEdits = @{\n    1  = @(\"51\u00f74=\", \"42\u00f76=\", \"20\u00f73=\", \"83\u00f72=\", \"15\u00f73=\", \"98\u00f75=\", \"10\u00f78=\", \"48\u00f75=\", \"88\u00f78=\", \"61\u00f77=\")\n    5  = @(\"34\u00f79=\", \"71\u00f79=\", \"14\u00f78=\", \"66\u00f76=\", \"48\u00f74=\", \"91\u00f72=\", \"40\u00f73=\", \"20\u00f77=\", \"31\u00f78=\", \"59\u00f73=\")\n    9  = @(\"98\u00f73=\", \"95\u00f74=\", \"30\u00f74=\", \"73\u00f76=\", \"31\u00f76=\", \"36\u00f77=\", \"97\u00f79=\", \"56\u00f78=\", \"13\u00f73=\", \"28\u00f77=\")\n    13 = @(\"82\u00f74=\", \"20\u00f77=\", \"24\u00f74=\", \"84\u00f76=\", \"64\u00f72=\", \"48\u00f74=\", \"33\u00f78=\", \"53\u00f76=\", \"88\u00f75=\", \"27\u00f74=\")\n    17 = @(\"72\u00f74=\", \"17\u00f77=\", \"15\u00f73=\", \"76\u00f79=\", \"67\u00f76=\", \"29\u00f73=\", \"37\u00f77=\", \"91\u00f75=\", \"60\u00f72=\", \"59\u00f75=\")\n}\n\nforeach ($rowIndex in $rowEdits.Keys) {\n    $pairs = $rowEdits[$rowIndex]\n    for ($col = 1; $col -le 5; $col++) {\n        $oldText = $pairs[($col - 1) * 2]\n        $newText = $pairs[($col - 1) * 2 + 1]\n        $cellRange = $table.Cell($rowIndex, $col).Range\n        $cellText = $cellRange.Text.TrimEnd([char]13, [char]7)\n        if ($cellText -eq $oldText) {\n            $cellRange.Text = $newText\n        }\n    }\n}\n"}
